$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that look numeric stay text (matches source formatting)
$ws.Range("D4:D6").NumberFormat = "@"
$ws.Range("D8:D10").NumberFormat = "@"
$ws.Range("D12:D14").NumberFormat = "@"
$ws.Range("D19:D25").NumberFormat = "@"
$ws.Range("D27:D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.699.29'
$ws.Range("E2").Value = '  +8.61%  '

$ws.Range("D3").Value = '3.480.20'
$ws.Range("E3").Value = '  +12.49%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '188.32'
$ws.Range("E5").Value = '  +13.34%  '

$ws.Range("D6").Value = '547.56'
$ws.Range("E6").Value = '  +8.18%  '

$ws.Range("D7").Value = '3.465.73'
$ws.Range("E7").Value = '  +12.30%  '

$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +4.39%  '

$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.22%  '

$ws.Range("D10").Value = '0.628'
$ws.Range("E10").Value = '  +8.88%  '

$ws.Range("E11").Value = '  +20.09%  '

$ws.Range("D12").Value = '54.53'
$ws.Range("E12").Value = '  +7.23%  '

$ws.Range("D13").Value = '0.0000267'
$ws.Range("E13").Value = '  +11.12%  '

$ws.Range("D14").Value = '9.32'
$ws.Range("E14").Value = '  +7.59%  '

$ws.Range("D15").Value = '4.020.44'
$ws.Range("E15").Value = '  +12.12%  '

$ws.Range("D16").Value = '3.467.10'
$ws.Range("E16").Value = '  +12.12%  '

$ws.Range("E17").Value = '  +8.15%  '

$ws.Range("D18").Value = '66.613.55'
$ws.Range("E18").Value = '  +8.83%  '

$ws.Range("D19").Value = '18.04'
$ws.Range("E19").Value = '  +9.36%  '

$ws.Range("D20").Value = '11.73'
$ws.Range("E20").Value = '  +11.82%  '

$ws.Range("D21").Value = '0.989'
$ws.Range("E21").Value = '  +6.84%  '

$ws.Range("D22").Value = '422.56'
$ws.Range("E22").Value = '  +19.66%  '

$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '3.88'
$ws.Range("E23").Value = '  +8.35%  '

$ws.Range("D24").Value = '84.24'
$ws.Range("E24").Value = '  +7.97%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '4.16'
$ws.Range("E25").Value = '  +9.24%  '

$ws.Range("E26").Value = '  +3.28%  '

$ws.Range("D27").Value = '2.89'
$ws.Range("E27").Value = '  +14.96%  '

$ws.Range("D28").Value = '6.14'
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").Value = '11.87'
$ws.Range("E29").Value = '  +10.87%  '

$ws.Range("D30").Value = '8.78'
$ws.Range("E30").Value = '  +12.64%  '

$ws.Range("D31").Value = '29.97'
$ws.Range("E31").Value = '  +10.48%  '

$ws.Range("D32").Value = '650.92'
$ws.Range("E32").Value = '  +4.45%  '

$ws.Range("D33").Value = '6.62'
$ws.Range("E33").Value = '  +7.67%  '

$ws.Range("D34").Value = '11.65'
$ws.Range("E34").Value = '  +7.13%  '

$ws.Range("D35").Value = '0.110'
$ws.Range("E35").Value = '  +9.70%  '

$ws.Range("D36").Value = '59.01'
$ws.Range("E36").Value = '  +6.20%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0813'
$ws.Range("E37").Value = '  +23.95%  '

$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").Value = '38.21'
$ws.Range("E38").Value = '  +9.66%  '

$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("D40").Value = '0.388'
$ws.Range("E40").Value = '  +7.30%  '

$ws.Range("E41").Value = '  +16.83%  '

$ws.Range("D42").Value = '3.33'
$ws.Range("E42").Value = '  +18.41%  '

$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.994.54'
$ws.Range("E44").Value = '  +8.20%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.43'
$ws.Range("E45").Value = '  +19.53%  '

$ws.Range("D46").Value = '2.62'
$ws.Range("E46").Value = '  +7.91%  '

$ws.Range("D47").Value = '2.88'
$ws.Range("E47").Value = '  +18.39%  '

$ws.Range("D48").Value = '0.0414'
$ws.Range("E48").Value = '  +10.85%  '

$ws.Range("D49").Value = '2.69'
$ws.Range("E49").Value = '  +3.86%  '

$ws.Range("D50").Value = '8.77'
$ws.Range("E50").Value = '  +20.44%  '

$ws.Range("D51").Value = '0.130'
$ws.Range("E51").Value = '  +8.78%  '
